$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) Paragraph 1 ("2) W" + bookmarkStart/End "_GoBack" + "e decided...")
#    -> drop the stray _GoBack bookmark, keep the text identical.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End
$p1Range = $d.Range($p1Start, $p1End)
$null = $p1Range.Delete()

$p1xml = "<w:p $wns>" `
    + "<w:r><w:t>2) W</w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`">e decided to use this wire-frame design because </w:t></w:r>" `
    + "<w:r><w:t>they were simple and not cluttered but still presented all of the necessary functionality.</w:t></w:r>" `
    + "</w:p>"
$null = $d.Range($p1Start, $p1Start).InsertXML($p1xml)

# ---------------------------------------------------------------------
# 2) Paragraphs "A)" / "Priming..." / "B)" / "i)..." / "ii)..." get
#    reshuffled & extended into: "A)+Priming", "B) ", "i)...supposed to...",
#    "ii)...(extended)", "iii)...(new, with the _GoBack bookmark moved here)".
# ---------------------------------------------------------------------
$pA = $d.Paragraphs(5)
$pLast = $d.Paragraphs(9)
$rangeStart = $pA.Range.Start
$rangeEnd = $pLast.Range.End
$null = $d.Range($rangeStart, $rangeEnd).Delete()

$frag = ""
$frag += "<w:p $wns>" `
    + "<w:r><w:t xml:space=`"preserve`">A) </w:t></w:r>" `
    + "<w:r><w:t>Priming: if you show a person an image, then a memory</w:t></w:r>" `
    + "<w:r><w:t>/association</w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`"> will be drawn up from their mind about it.</w:t></w:r>" `
    + "</w:p>"
$frag += "<w:p $wns>" `
    + "<w:r><w:t>B)</w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" `
    + "</w:p>"
$frag += "<w:p $wns>" `
    + "<w:proofErr w:type=`"spellStart`"/>" `
    + "<w:r><w:t>i</w:t></w:r>" `
    + "<w:proofErr w:type=`"spellEnd`"/>" `
    + "<w:r><w:t xml:space=`"preserve`">) </w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`">The principle could be applied </w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`">by keeping images professional, since the site we are making is </w:t></w:r>" `
    + "<w:r><w:t>supposed to</w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`"> be a business website.</w:t></w:r>" `
    + "</w:p>"
$frag += "<w:p $wns>" `
    + "<w:r><w:t>ii)</w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" `
    + "<w:r><w:t>The principle can lead to the success of the site</w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`">by making people who use the site feel a certain way. As in the previous point, we could make our website professional looking and have people </w:t></w:r>" `
    + "<w:proofErr w:type=`"gramStart`"/>" `
    + "<w:r><w:t>think</w:t></w:r>" `
    + "<w:proofErr w:type=`"gramEnd`"/>" `
    + "<w:r><w:t xml:space=`"preserve`"> of that.</w:t></w:r>" `
    + "</w:p>"
$frag += "<w:p $wns>" `
    + "<w:r><w:t>iii)</w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`"> This principle could be </w:t></w:r>" `
    + "<w:r><w:t xml:space=`"preserve`">also </w:t></w:r>" `
    + "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" `
    + "<w:bookmarkEnd w:id=`"0`"/>" `
    + "<w:r><w:t>applied in subtle ways such as putting a lock by the password to help people make a strong password.</w:t></w:r>" `
    + "</w:p>"

$null = $d.Range($rangeStart, $rangeStart).InsertXML($frag)
